# Generate Report for Handback
# Updates the localization-status report after a handback:
#  - Overview sheet: Status for zh-cn/de-de rows moves from "Ready for
#    handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: fill in "Latest Target File" and
#    "Latest Handback File" links/names, and stamp "Latest Handback
#    DateTime" for each row

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = "a.md"
$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-23 22:34:34"

$zh.Range("I3").Value = "a.md"
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-23 22:34:34"

# Rebuild the hyperlinks on this sheet in row order (A2, I2, A3, I3) so
# the relationship ids line up the way Excel assigns them.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/b.md", $null, $null, "b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")

# ---- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = "a.md"
$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K2").Value = "2016-08-23 22:34:42"

$de.Range("I3").Value = "a.md"
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K3").Value = "2016-08-23 22:34:42"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/b.md", $null, $null, "b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md", $null, $null, "a.md")

# ---- Column widths: widen columns whose content grew -----------------
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zh.Columns.Item(3).AutoFit()
$zh.Columns.Item(10).AutoFit()
$de.Columns.Item(3).AutoFit()
$de.Columns.Item(10).AutoFit()
